$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 6814327
$ws.Cells.Item(2, 5).Value = 'NS Mura'
$ws.Cells.Item(2, 6).Value = 'NK Domzale'
$ws.Cells.Item(2, 8).Value = 3
$ws.Cells.Item(2, 9).Value = 'A'
$ws.Cells.Item(2, 10).Value = 2
$ws.Cells.Item(2, 11).Value = 3.3
$ws.Cells.Item(2, 12).Value = 3.4
$ws.Cells.Item(2, 13).Value = 1.909
$ws.Cells.Item(2, 14).Value = 3.4
$ws.Cells.Item(2, 16).Value = -0.5
$ws.Cells.Item(2, 17).Value = 1.95
$ws.Cells.Item(2, 18).Value = 1.85
$ws.Cells.Item(2, 19).Value = 2.5
$ws.Cells.Item(2, 20).Value = 1.9
$ws.Cells.Item(2, 21).Value = 1.9
$ws.Cells.Item(2, 22).Value = -1
$ws.Cells.Item(2, 24).Value = 2.75
$ws.Cells.Item(2, 25).Value = -1
$ws.Cells.Item(2, 26).Value = 0.8500000000000001
$ws.Cells.Item(2, 27).Value = 0.8999999999999999
$ws.Cells.Item(2, 28).Value = -1
$ws.Cells.Item(3, 2).Value = 6816473
$ws.Cells.Item(3, 5).Value = 'NK Bravo'
$ws.Cells.Item(3, 6).Value = 'NK Rogaska'
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 'H'
$ws.Cells.Item(3, 10).Value = 1.8
$ws.Cells.Item(3, 11).Value = 3.5
$ws.Cells.Item(3, 12).Value = 4
$ws.Cells.Item(3, 13).Value = 2.05
$ws.Cells.Item(3, 14).Value = 3
$ws.Cells.Item(3, 16).Value = -0.25
$ws.Cells.Item(3, 17).Value = 1.75
$ws.Cells.Item(3, 18).Value = 2.05
$ws.Cells.Item(3, 19).Value = 2.25
$ws.Cells.Item(3, 20).Value = 1.95
$ws.Cells.Item(3, 21).Value = 1.85
$ws.Cells.Item(3, 22).Value = 1.05
$ws.Cells.Item(3, 24).Value = -1
$ws.Cells.Item(3, 25).Value = 0.75
$ws.Cells.Item(3, 26).Value = -1
$ws.Cells.Item(3, 27).Value = -0.5
$ws.Cells.Item(3, 28).Value = 0.425
$ws.Cells.Item(9, 2).Value = 6814328
$ws.Cells.Item(9, 5).Value = 'NK Domzale'
$ws.Cells.Item(9, 6).Value = 'NK Bravo'
$ws.Cells.Item(9, 8).Value = 1
$ws.Cells.Item(9, 9).Value = 'D'
$ws.Cells.Item(9, 10).Value = 2.35
$ws.Cells.Item(9, 11).Value = 3.1
$ws.Cells.Item(9, 12).Value = 2.9
$ws.Cells.Item(9, 13).Value = 2.15
$ws.Cells.Item(9, 14).Value = 3.1
$ws.Cells.Item(9, 15).Value = 3.3
$ws.Cells.Item(9, 16).Value = -0.25
$ws.Cells.Item(9, 17).Value = 1.925
$ws.Cells.Item(9, 18).Value = 1.875
$ws.Cells.Item(9, 19).Value = 2.25
$ws.Cells.Item(9, 20).Value = 1.95
$ws.Cells.Item(9, 21).Value = 1.85
$ws.Cells.Item(9, 22).Value = -1
$ws.Cells.Item(9, 23).Value = 2.1
$ws.Cells.Item(9, 26).Value = 0.4375
$ws.Cells.Item(9, 27).Value = -0.5
$ws.Cells.Item(9, 28).Value = 0.425
$ws.Cells.Item(10, 2).Value = 6814330
$ws.Cells.Item(10, 5).Value = 'NK Maribor'
$ws.Cells.Item(10, 6).Value = 'NK Aluminij'
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 'H'
$ws.Cells.Item(10, 10).Value = 1.363
$ws.Cells.Item(10, 11).Value = 4.5
$ws.Cells.Item(10, 12).Value = 7
$ws.Cells.Item(10, 13).Value = 1.4
$ws.Cells.Item(10, 14).Value = 4.5
$ws.Cells.Item(10, 15).Value = 7
$ws.Cells.Item(10, 16).Value = -1.25
$ws.Cells.Item(10, 17).Value = 1.85
$ws.Cells.Item(10, 18).Value = 1.95
$ws.Cells.Item(10, 19).Value = 2.75
$ws.Cells.Item(10, 20).Value = 1.8
$ws.Cells.Item(10, 21).Value = 2
$ws.Cells.Item(10, 22).Value = 0.3999999999999999
$ws.Cells.Item(10, 23).Value = -1
$ws.Cells.Item(10, 26).Value = 0.475
$ws.Cells.Item(10, 27).Value = -1
$ws.Cells.Item(10, 28).Value = 1
$ws.Cells.Item(159, 17).Value = 1.9
$ws.Cells.Item(159, 18).Value = 1.9
$ws.Cells.Item(159, 20).Value = 1.9
$ws.Cells.Item(159, 21).Value = 1.9
$ws.Cells.Item(160, 13).Value = 4.75
$ws.Cells.Item(160, 14).Value = 3.75
$ws.Cells.Item(160, 15).Value = 1.615
$ws.Cells.Item(160, 17).Value = 2.025
$ws.Cells.Item(160, 18).Value = 1.775
$ws.Cells.Item(160, 20).Value = 1.95
$ws.Cells.Item(160, 21).Value = 1.85
$ws.Cells.Item(161, 13).Value = 1.909
$ws.Cells.Item(161, 15).Value = 3.8
$ws.Cells.Item(161, 17).Value = 1.925
$ws.Cells.Item(161, 18).Value = 1.875
$ws.Cells.Item(162, 20).Value = 1.95
$ws.Cells.Item(162, 21).Value = 1.85
$ws.Cells.Item(163, 13).Value = 1.95
$ws.Cells.Item(163, 17).Value = 2
$ws.Cells.Item(163, 18).Value = 1.8
$ws.Cells.Item(163, 20).Value = 1.9
$ws.Cells.Item(163, 21).Value = 1.9
